# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date", Priority, and handoff/handback
# datetime values that get refreshed each time the handback status report is
# regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" column (G) on the Overview sheet.
# G2 and G3 previously shared the value "2016-09-06 06:18:16"; both rows
# move to the new generation timestamp together.
$wsOverview.Range("G2").Value = "2016-09-06 06:18:59"
$wsOverview.Range("G3").Value = "2016-09-06 06:18:59"

# "Priority" column (E) changes from "ht" (human translation) to
# "mt" (machine translation) for both language sheets / both rows.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns refresh to new timestamps.
$wsZhCn.Range("H2").Value = "2016-09-06 06:18:55"
$wsZhCn.Range("H3").Value = "2016-09-06 06:18:55"
$wsZhCn.Range("K2").Value = "2016-09-06 06:19:19"
$wsZhCn.Range("K3").Value = "2016-09-06 06:19:19"

# de-de sheet: "Correspond Handoff Datetime" (H) refreshes to the new
# generation timestamp (shared with Overview!G2:G3).
$wsDeDe.Range("H2").Value = "2016-09-06 06:18:59"
$wsDeDe.Range("H3").Value = "2016-09-06 06:18:59"

# de-de sheet: "Correspond Handback DateTime" (K) refreshes to a new
# timestamp.
$wsDeDe.Range("K2").Value = "2016-09-06 06:19:26"
$wsDeDe.Range("K3").Value = "2016-09-06 06:19:26"
